$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 741, shifting existing rows 741-846 down to 742-847
$ws.Rows.Item(741).Insert()

# Fill in the new row 741 with data
$ws.Cells.Item(741, 1).Value = 10
$ws.Cells.Item(741, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(741, 3).Value = "La Araucanía"
$ws.Cells.Item(741, 4).Value = 45212
$ws.Cells.Item(741, 5).Value = 9
$ws.Cells.Item(741, 6).Value = 100112027
$ws.Cells.Item(741, 7).Value = "Melón"
$ws.Cells.Item(741, 8).Value = "Tuna"
$ws.Cells.Item(741, 9).Value = "Primera"
$ws.Cells.Item(741, 10).Value = 150
$ws.Cells.Item(741, 11).Value = 1900
$ws.Cells.Item(741, 12).Value = 1900
$ws.Cells.Item(741, 13).Value = 1900
$ws.Cells.Item(741, 14).Value = "$/unidad"
$ws.Cells.Item(741, 15).Value = "Perú"
$ws.Cells.Item(741, 16).Value = 1900
$ws.Cells.Item(741, 17).Value = 1
$ws.Cells.Item(741, 18).Value = "Hortaliza"
